$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply custom-accuracy rounding (2 decimal places) to row 5 measurement values
$ws.Range("B5").Value = 10.88
$ws.Range("C5").Value = 8.18
$ws.Range("D5").Value = 0.25
$ws.Range("E5").Value = 24.03
$ws.Range("F5").Value = 19.95
$ws.Range("G5").Value = 8.36
$ws.Range("H5").Value = 33.74
$ws.Range("I5").Value = 13.18
$ws.Range("J5").Value = 6.01
$ws.Range("K5").Value = 9.01
$ws.Range("L5").Value = 10.09
$ws.Range("M5").Value = 10.3
$ws.Range("N5").Value = 3.09
$ws.Range("O5").Value = 8.48
$ws.Range("P5").Value = 12.33
$ws.Range("Q5").Value = 7.09
$ws.Range("R5").Value = 0.13
$ws.Range("S5").Value = 0.5
$ws.Range("T5").Value = 124.07
$ws.Range("U5").Value = 23.85
$ws.Range("V5").Value = 7.82
$ws.Range("W5").Value = 16.26
$ws.Range("X5").Value = 8.42
$ws.Range("Y5").Value = 1.14
$ws.Range("Z5").Value = 16.38
$ws.Range("AA5").Value = 7
$ws.Range("AB5").Value = 6.39
$ws.Range("AC5").Value = 7.38
$ws.Range("AD5").Value = 10.53
$ws.Range("AE5").Value = 0.29
$ws.Range("AF5").Value = 31
$ws.Range("AG5").Value = 4.71
$ws.Range("AH5").Value = 9.81

# Remove the now-unused last data row (row 6) - dataset trimmed to 1000 rows upstream
$ws.Rows("6:6").Delete()
